$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rRating = "Children under 17 may not attend R-rated movies unaccompanied by a parent or adult guardian. Children ages 6 and under are not allowed at R-rated movies after 6pm.`r`nR"

$movies = @(
    ,("TAYLOR SWIFT | THE ERAS TOUR", "2 HR 48 MIN", "PG13", "Released Oct 13, 2023")
    ,("TAYLOR SWIFT | THE ERAS TOUR: Private Theatre Rental for up to 40 Guests", "2 HR 48 MIN", "PG13", "Released Oct 13, 2023")
    ,("The Exorcist: Believer", "1 HR 51 MIN", $rRating, "Released Oct 6, 2023")
    ,("Paw Patrol: The Mighty Movie", "1 HR 28 MIN", "PG", "Released Sep 29, 2023")
    ,("Saw X", "1 HR 58 MIN", $rRating, "Released Sep 29, 2023")
    ,("The Nun II", "1 HR 50 MIN", $rRating, "Released Sep 8, 2023")
    ,("The Creator", "2 HR 13 MIN", "PG13", "Released Sep 29, 2023")
    ,("Dumb Money", "1 HR 45 MIN", $rRating, "Released Sep 15, 2023")
    ,("A Haunting in Venice", "1 HR 43 MIN", "PG13", "Released Sep 15, 2023")
    ,("The Equalizer 3", "1 HR 49 MIN", $rRating, "Released Sep 1, 2023")
    ,("Barbie", "1 HR 54 MIN", "PG13", "Released Jul 21, 2023")
    ,("EXPEND4BLES", "1 HR 43 MIN", $rRating, "Released Sep 22, 2023")
    ,("My Big Fat Greek Wedding 3", "1 HR 32 MIN", "PG13", "Released Sep 8, 2023")
    ,("Oppenheimer", "3 HR", $rRating, "Released Jul 21, 2023")
    ,("Hocus Pocus 30th Anniversary", "1 HR 36 MIN", "PG", "Released Oct 6, 2023")
    ,("Under The Skin", "1 HR 50 MIN", $rRating, "Opening Oct 18, 2023")
    ,("Beetlejuice 35th Anniversary", "1 HR 32 MIN", "PG", "Opening Oct 18, 2023")
    ,("The Blind", "1 HR 58 MIN", "PG13", "Released Sep 28, 2023")
    ,("Still Playin' Possum: Music and Memories of George Jones", "1 HR 45 MIN", "NR", "Opening Oct 17, 2023")
    ,("When Evil Lurks", "1 HR 40 MIN", "NR", "Released Oct 6, 2023")
    ,("Stop Making Sense", "1 HR 28 MIN", "PG", "Released Sep 22, 2023")
    ,("Strange Way of Life - With Added Bonus Attraction: The Human Voice", "1 HR 1 MIN", $rRating, "Released Oct 6, 2023")
    ,("She Came to Me", "1 HR 42 MIN", "NR", "Released Oct 6, 2023")
    ,("Shelter in Solitude", "1 HR 33 MIN", "NR", "Released Oct 6, 2023")
    ,("Private Theatre Rental", $null, $null, $null)
)

for ($i = 0; $i -lt $movies.Count; $i++) {
    $row = $i + 1
    $data = $movies[$i]
    for ($c = 0; $c -lt 4; $c++) {
        $val = $data[$c]
        if ($null -ne $val) {
            $ws.Cells.Item($row, $c + 1).Value = $val
        } else {
            $ws.Cells.Item($row, $c + 1).Value = $null
        }
    }
}
